$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# This is a weekly data refresh: the "Fecha" (and associated Volumen /
# Precio / Precio-$-por-Kg) values roll forward from one reporting row to
# another. The net effect on the data is a rotation through these row
# groups (each row keeps its Calidad/other descriptive columns, only the
# date-dependent measurement columns D, M, N, O, P, S move):
#   2 -> 15 -> 4 -> 13 -> 2      (row 2 ends up holding what row 13 had,
#                                  row 13 ends up holding what row 4 had,
#                                  row 4 ends up holding what row 15 had,
#                                  row 15 ends up holding what row 2 had)
#   3 -> 16 -> 5 -> 14 -> 3
#   11 -> 17 -> 11
#   12 -> 18 -> 12
#
# i.e. destination row gets the values previously held by the "source" row
# listed next to it below.

$cols = @("D", "M", "N", "O", "P", "S")

# Each inner array lists rows in the order data flows: the first row
# receives the values the last row had (cyclic rotation), i.e.
# destination[i] <= source values currently sitting in cycle[i+1].
$cycles = @(
    @(2, 13, 4, 15),
    @(3, 14, 5, 16),
    @(11, 17),
    @(12, 18)
)

foreach ($cycle in $cycles) {
    foreach ($col in $cols) {
        # Snapshot current values for every row in the cycle first.
        $snapshot = @{}
        foreach ($row in $cycle) {
            $snapshot[$row] = $ws.Range("$col$row").Value2
        }

        # destRow receives the value that used to belong to the next row
        # in the cycle (wrapping around).
        for ($i = 0; $i -lt $cycle.Length; $i++) {
            $destRow = $cycle[$i]
            $srcRow = $cycle[($i + 1) % $cycle.Length]
            $ws.Range("$col$destRow").Value = $snapshot[$srcRow]
        }
    }
}
